$d = $word.ActiveDocument

# 1. Update years of experience in PROFESSIONAL SUMMARY
$d.Content.Find.Execute("Distinguished Polling, Research & Redistricting Professional with 21 years of expertise", $true, $false, $false, $false, $false, $true, 1, $false, "Distinguished Polling, Research & Redistricting Professional with 15+ years of expertise", 2) | Out-Null

# 2. Update FLEEM web application bullet (Progressive Change Campaign Committee)
$d.Content.Find.Execute("• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls", $true, $false, $false, $false, $false, $true, 1, $false, "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys", 2) | Out-Null

# 3. Expand Salsa Labs bullet list
$old3 = "• Developed software solutions for political campaigns and advocacy groups^p• Built web applications for voter engagement and campaign management^p• Integrated third-party APIs and data sources for campaign tools^p• Collaborated with political strategists to translate requirements into technical solutions"
$new3 = "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously^p• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers^p• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill^p• Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs^p• Collaborated with political strategists to translate geospatial requirements into technical solutions^p• Handled billions of records with millions of columns in high-performance CRM system"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# 4. Expand Praxis Project bullet list
$old4 = "• Integrated technology solutions within organizational frameworks for social justice organizations^p• Developed data management systems for community organizing efforts^p• Provided technical training and support to nonprofit staff^p• Built custom applications for community engagement and advocacy"
$new4 = "• Led technology operations for multi-million dollar organization while assisting in search for full-time CTO^p• Directed all technology decisions and practices for massive multinational non-governmental organization^p• Developed comprehensive frameworks for internal and external technology audits^p• Led training initiatives for beneficiaries on spatial and Census data analysis for public health research^p• Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL^p• Managed technology infrastructure supporting community health initiatives across multiple countries^p• Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation"
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# 5. Add new bullet after Lake Research Partners last bullet
$old5 = "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding"
$new5 = "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding^p• Trained staff on building Python tooling for report generation and analysis"
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

# 6. Add new bullet after Feldman Group last bullet
$old6 = "• Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL"
$new6 = "• Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL^p• Trained staff on PHP/MySQL for data analysis and reporting systems"
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

Write-Output "All replacements complete"
